$wb = $excel.ActiveWorkbook

$wsTests = $wb.Worksheets.Item("Tests")
$wsResult = $wb.Worksheets.Item("Result")

# Add a new row (10) to the Tests sheet with the KillAllProcesses workflow file
$wsTests.Range("A10").Value = "Framework\KillAllProcesses.xaml"
$wsTests.Range("B10").Value = "Success"

# Add a new row (10) to the Result sheet with the KillAllProcesses workflow file
$wsResult.Range("A10").Value = "Framework\KillAllProcesses.xaml"
$wsResult.Range("B10").Value = "Success"

# Update selections to match target state
$wsTests.Range("A30").Select()
$wsResult.Range("D18").Select()
